# Update "want to go" counts (column F) in several worksheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    "展览" = @{
        4 = 3756; 5 = 3756; 6 = 292; 7 = 5294; 8 = 586; 9 = 423; 10 = 227;
        11 = 1057; 13 = 138; 15 = 729; 19 = 176; 22 = 6052; 23 = 6052;
        25 = 43; 27 = 6818; 30 = 3256; 31 = 368; 32 = 749; 33 = 4459;
        34 = 325; 36 = 151; 37 = 1148; 39 = 28; 41 = 922; 42 = 1135; 43 = 2059
    }
    "演出" = @{
        3 = 31; 5 = 66
    }
    "本地生活" = @{
        2 = 231; 3 = 1153
    }
    "全部类型" = @{
        2 = 231; 4 = 1153; 7 = 3756; 8 = 3756; 9 = 292; 10 = 5294; 11 = 586;
        12 = 423; 13 = 227; 14 = 1057; 16 = 138; 18 = 729; 23 = 176; 26 = 6052;
        30 = 6818; 33 = 3256; 34 = 368; 35 = 749; 36 = 4459; 37 = 325; 38 = 31;
        40 = 151; 41 = 1148; 43 = 28; 45 = 922; 46 = 1135; 48 = 2059; 50 = 66
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Cells.Item($row, 6).Value = $rows[$row]
    }
}
